$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update subtitle (row 2) - "October 2016" -> "November 2016"
$ws.Range("A2").Value = "by State, by Sector, Year-to-Date through November 2016 and 2015 (Thousand Megawatthours)"

# Update column headers in row 6 (all cells sharing "October 2016 YTD" / "October 2015 YTD")
foreach ($cellRef in @("B6","E6","G6","I6","K6")) {
    $ws.Range($cellRef).Value = "November 2016 YTD"
}
foreach ($cellRef in @("C6","F6","H6","J6","L6")) {
    $ws.Range($cellRef).Value = "November 2015 YTD"
}

# Middle Atlantic
$ws.Range("B14").Value = 672
$ws.Range("C14").Value = 593
$ws.Range("D14").Value = 0.133
$ws.Range("K14").Value = 671
$ws.Range("L14").Value = 593
# New Jersey
$ws.Range("B15").Value = 210
$ws.Range("C15").Value = 204
$ws.Range("D15").Value = 0.03
$ws.Range("K15").Value = 210
$ws.Range("L15").Value = 204
# Pennsylvania
$ws.Range("B17").Value = 461
$ws.Range("C17").Value = 389
$ws.Range("D17").Value = 0.187
$ws.Range("K17").Value = 460
$ws.Range("L17").Value = 389
# East North Central
$ws.Range("B18").Value = 4550
$ws.Range("C18").Value = 4383
$ws.Range("D18").Value = 0.038
$ws.Range("E18").Value = 142
$ws.Range("F18").Value = 197
$ws.Range("G18").Value = 2018
$ws.Range("H18").Value = 1685
$ws.Range("K18").Value = 2390
$ws.Range("L18").Value = 2500
# Illinois
$ws.Range("B19").Value = 266
$ws.Range("C19").Value = 236
$ws.Range("D19").Value = 0.127
$ws.Range("H19").Value = 2
$ws.Range("K19").Value = 261
$ws.Range("L19").Value = 234
# Indiana
$ws.Range("B20").Value = 1954
$ws.Range("C20").Value = 2103
$ws.Range("D20").Value = -0.071
$ws.Range("F20").Value = 19
$ws.Range("K20").Value = 1935
$ws.Range("L20").Value = 2084
# Michigan
$ws.Range("B21").Value = 1478
$ws.Range("D21").Value = 0.257
$ws.Range("E21").Value = 123
$ws.Range("G21").Value = 1355
# Ohio
$ws.Range("B22").Value = 852
$ws.Range("C22").Value = 868
$ws.Range("D22").Value = -0.019
$ws.Range("G22").Value = 658
$ws.Range("H22").Value = 686
$ws.Range("K22").Value = 194
$ws.Range("L22").Value = 182
# West North Central
$ws.Range("B24").Value = 43
$ws.Range("C24").Value = 37
$ws.Range("D24").Value = 0.163
$ws.Range("K24").Value = 43
$ws.Range("L24").Value = 37
# North Dakota
$ws.Range("B30").Value = 43
$ws.Range("C30").Value = 37
$ws.Range("D30").Value = 0.163
$ws.Range("K30").Value = 43
$ws.Range("L30").Value = 37
# South Atlantic
$ws.Range("B32").Value = 280
$ws.Range("C32").Value = 245
$ws.Range("D32").Value = 0.146
$ws.Range("K32").Value = 280
$ws.Range("L32").Value = 245
# Delaware
$ws.Range("B33").Value = 253
$ws.Range("C33").Value = 213
$ws.Range("D33").Value = 0.19
$ws.Range("K33").Value = 253
$ws.Range("L33").Value = 213
# Florida
$ws.Range("B35").Value = 5
$ws.Range("C35").Value = 5
$ws.Range("D35").Value = -0.006
$ws.Range("K35").Value = 5
$ws.Range("L35").Value = 5
# West Virginia
$ws.Range("B41").Value = 22
$ws.Range("C41").Value = 27
$ws.Range("D41").Value = -0.175
$ws.Range("K41").Value = 22
$ws.Range("L41").Value = 27
# East South Central
$ws.Range("B42").Value = 37
$ws.Range("C42").Value = 47
$ws.Range("D42").Value = -0.204
$ws.Range("K42").Value = 37
$ws.Range("L42").Value = 47
# Alabama
$ws.Range("B43").Value = 23
$ws.Range("D43").Value = -0.35
$ws.Range("K43").Value = 23
# Tennessee
$ws.Range("B46").Value = 14
$ws.Range("C46").Value = 11
$ws.Range("D46").Value = 0.271
$ws.Range("K46").Value = 14
$ws.Range("L46").Value = 11
# West South Central
$ws.Range("B47").Value = 4330
$ws.Range("C47").Value = 4496
$ws.Range("D47").Value = -0.037
$ws.Range("G47").Value = 1206
$ws.Range("H47").Value = 1149
$ws.Range("K47").Value = 3124
$ws.Range("L47").Value = 3347
# Louisiana
$ws.Range("B49").Value = 1922
$ws.Range("C49").Value = 2181
$ws.Range("D49").Value = -0.119
$ws.Range("K49").Value = 1922
$ws.Range("L49").Value = 2181
# Texas
$ws.Range("B51").Value = 2408
$ws.Range("C51").Value = 2315
$ws.Range("D51").Value = 0.04
$ws.Range("G51").Value = 1206
$ws.Range("H51").Value = 1149
$ws.Range("K51").Value = 1202
$ws.Range("L51").Value = 1166
# Mountain
$ws.Range("B52").Value = 335
$ws.Range("C52").Value = 396
$ws.Range("D52").Value = -0.153
$ws.Range("G52").Value = 8
$ws.Range("H52").Value = 21
$ws.Range("K52").Value = 327
$ws.Range("L52").Value = 375
# Montana
$ws.Range("B56").Value = 7
$ws.Range("C56").Value = 16
$ws.Range("D56").Value = -0.537
$ws.Range("G56").Value = 7
$ws.Range("H56").Value = 16
# Nevada
$ws.Range("D57").Value = -0.862
# Utah
$ws.Range("C59").Value = 8
$ws.Range("L59").Value = 8
# Wyoming
$ws.Range("B60").Value = 321
$ws.Range("C60").Value = 368
$ws.Range("D60").Value = -0.127
$ws.Range("K60").Value = 321
$ws.Range("L60").Value = 368
# Pacific Contiguous
$ws.Range("B61").Value = 1700
$ws.Range("C61").Value = 1764
$ws.Range("D61").Value = -0.036
$ws.Range("G61").Value = 364
$ws.Range("H61").Value = 359
$ws.Range("K61").Value = 1336
$ws.Range("L61").Value = 1405
# California
$ws.Range("B62").Value = 1336
$ws.Range("C62").Value = 1405
$ws.Range("D62").Value = -0.049
$ws.Range("K62").Value = 1336
$ws.Range("L62").Value = 1405
# Washington
$ws.Range("B64").Value = 364
$ws.Range("C64").Value = 359
$ws.Range("D64").Value = 0.014
$ws.Range("G64").Value = 364
$ws.Range("H64").Value = 359
# Pacific Noncontiguous
$ws.Range("B65").Value = 39
$ws.Range("C65").Value = 46
$ws.Range("D65").Value = -0.142
$ws.Range("K65").Value = 39
$ws.Range("L65").Value = 46
# Hawaii
$ws.Range("B67").Value = 39
$ws.Range("C67").Value = 46
$ws.Range("D67").Value = -0.142
$ws.Range("K67").Value = 39
$ws.Range("L67").Value = 46
# U.S. Total
$ws.Range("B68").Value = 11987
$ws.Range("C68").Value = 12007
$ws.Range("D68").Value = -0.002
$ws.Range("E68").Value = 142
$ws.Range("F68").Value = 197
$ws.Range("G68").Value = 3598
$ws.Range("H68").Value = 3214
$ws.Range("K68").Value = 8248
$ws.Range("L68").Value = 8595

Write-Host "Edit complete"
